$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row 26: "Exp 28" experiment (Digits 0, 1, 9) ---
# Set the two text cells that introduce the first two new shared strings
# (so the shared-strings table grows in the same order as the source edit).
$ws.Range("A26").Value = "Exp 28"
$ws.Range("F26").Value = "Exp 28.png"

# --- New secondary header row 25 (G:K) introducing the "wrt 9" comparison columns ---
$ws.Range("G25").Value = "Micro 9"
$ws.Range("H25").Value = "Micro 1"
$ws.Range("I25").Value = "Micro 0"
$ws.Range("J25").Value = "Macro 1 wrt 9"
$ws.Range("K25").Value = "Macro 0 wrt 9"
$ws.Range("G25:K25").Style = $ws.Range("G1:K1").Style

# --- Remaining values for row 26 ---
$ws.Range("B26").Value = 0.4
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = "Local"
$ws.Range("E26").Value = -1
$ws.Range("G26").Value = 88.24
$ws.Range("H26").Value = 88.09
$ws.Range("I26").Value = 87.93
$ws.Range("J26").Value = 74.76
$ws.Range("K26").Value = 76.42

# Match the data-row styling used by the other experiment rows.
$ws.Range("A26:E26").Style = $ws.Range("A2:E2").Style
$ws.Range("G26:K26").Style = $ws.Range("G2:K2").Style

# Scroll/select so the new row is in view, matching the saved view state.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D26").Select()
